$wb = $excel.ActiveWorkbook

# --- Rename first sheet: Glider -> Moorings ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Moorings"

$ws2 = $wb.Worksheets.Item(2)

# --- Update Asset_Cal_Info (sheet2) Ref Des text values so the GL001 -> GL388 rename
#     propagates through the shared-string table ---
$ws2.Range("A2").Value = "CP05MOAS-GL388-01-ADCPAM000"
$ws2.Range("A3").Value = "CP05MOAS-GL388-01-ADCPAM000"
$ws2.Range("A4").Value = "CP05MOAS-GL388-01-ADCPAM000"
$ws2.Range("A5").Value = "CP05MOAS-GL388-01-ADCPAM000"
$ws2.Range("A7").Value = "CP05MOAS-GL388-02-FLORTM000"
$ws2.Range("A8").Value = "CP05MOAS-GL388-02-FLORTM000"
$ws2.Range("A9").Value = "CP05MOAS-GL388-02-FLORTM000"
$ws2.Range("A10").Value = "CP05MOAS-GL388-02-FLORTM000"
$ws2.Range("A12").Value = "CP05MOAS-GL388-03-CTDGVM000"
$ws2.Range("A14").Value = "CP05MOAS-GL388-04-DOSTAM000"
$ws2.Range("A16").Value = "CP05MOAS-GL388-05-PARADM000"
$ws2.Range("A18").Value = "CP05MOAS-GL388-00-ENG000000"

# --- Update Moorings (sheet1) row 2 values ---
$ws1.Range("E2").Value = 0.0625
$ws1.Range("I2").Value = 0

# --- Add new Latitude/Longitude decimal-degree formula columns ---
$ws1.Range("L2").Formula = '=((LEFT(G2,(FIND("°",G2,1)-1)))+(MID(G2,(FIND("°",G2,1)+1),(FIND("''",G2,1))-(FIND("°",G2,1)+1))/60))*(IF(RIGHT(G2,1)="N",1,-1))'
$ws1.Range("M2").Formula = '=((LEFT(H2,(FIND("°",H2,1)-1)))+(MID(H2,(FIND("°",H2,1)+1),(FIND("''",H2,1))-(FIND("°",H2,1)+1))/60))*(IF(RIGHT(H2,1)="E",1,-1))'

$newCellsRange = $ws1.Range("L2:M2")
$newCellsRange.Font.Size = 11
$newCellsRange.Font.Name = "Calibri"
$newCellsRange.Font.Color = 0
$newCellsRange.HorizontalAlignment = -4108
$newCellsRange.VerticalAlignment = -4108

# --- Selections to mirror the saved workbook views ---
$ws2.Range("E26").Select()
$ws1.Range("C20").Select()
